# Fruta / hortaliza, semanal
# The underlying data rows (2-20) were re-shuffled: each row's Fecha (D),
# Calidad (I), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) were swapped between
# rows according to a fixed permutation (row 21 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: new row -> old row whose D/I/J/K/L/M/P values it should receive.
$rowMap = @{
    2  = 12
    3  = 18
    4  = 15
    5  = 16
    6  = 20
    7  = 17
    8  = 2
    9  = 3
    10 = 11
    11 = 9
    12 = 10
    13 = 4
    14 = 5
    15 = 14
    16 = 19
    17 = 13
    18 = 7
    19 = 6
    20 = 8
    21 = 21
}

# Snapshot the original values for every affected row before writing
# anything, since several rows both give and receive data. Value2 is
# used (rather than Value) so plain dates/numbers/strings are captured
# instead of a live COM variant reference.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        I = $ws.Cells.Item($r, 9).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $src = $snapshot[$oldRow]

    $ws.Cells.Item($newRow, 4).Value2 = $src.D
    $ws.Cells.Item($newRow, 9).Value2 = $src.I
    $ws.Cells.Item($newRow, 10).Value2 = $src.J
    $ws.Cells.Item($newRow, 11).Value2 = $src.K
    $ws.Cells.Item($newRow, 12).Value2 = $src.L
    $ws.Cells.Item($newRow, 13).Value2 = $src.M
    $ws.Cells.Item($newRow, 16).Value2 = $src.P
}
